$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge the three split runs that together spell out the contractor-terms
#    URL ("...contractor/1e2u).") back into a single run. The three runs
#    already render as contiguous text, so a self-replace via Find/Replace
#    (same text in, same text out) is enough to coalesce them into one run.
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Text = "1e2u)."
$find.Replacement.Text = "1e2u)."
$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)

# ---------------------------------------------------------------------------
# 2) Remove the "{ContractorCo, Inc.}, / {a Delaware corporation}" entity
#    placeholder from the contractor signature block, leaving a blank
#    paragraph in its place (mirroring the already-blank line under the
#    client block).
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*{ContractorCo, Inc.}*") {
        $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:keepNext w:val="true"/><w:widowControl w:val="false"/><w:jc w:val="left"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/></w:rPr></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:widowControl w:val="false"/><w:jc w:val="left"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/></w:rPr></w:r></w:p>'
        $p.Range.InsertXML($xml)
        break
    }
}

# ---------------------------------------------------------------------------
# 3) Turn on "suppress automatic hyphens" for the Normal and LO-normal
#    paragraph styles.
# ---------------------------------------------------------------------------
$d.Styles("Normal").ParagraphFormat.Hyphenation = $false
$d.Styles("LO-normal").ParagraphFormat.Hyphenation = $false
